$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.756.81"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "1.894.29"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'0.7632"
$ws.Range("E5").Value = "  +4.08%  "
$ws.Range("D6").Value = "'240.40"
$ws.Range("E6").Value = "  -1.01%  "
$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "1.894.15"
$ws.Range("E8").Value = "  +0.97%  "
$ws.Range("D9").Value = "'0.3043"
$ws.Range("E9").Value = "  -1.90%  "
$ws.Range("D10").Value = "'25.33"
$ws.Range("E10").Value = "  -3.18%  "
$ws.Range("D11").Value = "'0.06807"
$ws.Range("E11").Value = "  -1.33%  "
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").Value = "1.890.18"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").Value = "'0.7350"
$ws.Range("E14").Value = "  -4.71%  "
$ws.Range("D15").Value = "'5.137"
$ws.Range("E15").Value = "  -1.60%  "
$ws.Range("D16").Value = "'90.70"
$ws.Range("E16").Value = "  -0.68%  "
$ws.Range("D17").Value = "29.778.32"
$ws.Range("D18").Value = "'13.81"
$ws.Range("E18").Value = "  -2.58%  "
$ws.Range("D19").Value = "'5.902"
$ws.Range("E19").Value = "  +2.74%  "
$ws.Range("D20").Value = "'241.59"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("D21").Value = "'0.000007673"
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'6.888"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").Value = "'166.15"
$ws.Range("D26").Value = "'9.170"
$ws.Range("E26").Value = "  -1.37%  "
$ws.Range("D27").Value = "'18.59"
$ws.Range("E27").Value = "  -1.37%  "
$ws.Range("D28").Value = "'0.1287"
$ws.Range("E28").Value = "  +1.65%  "
$ws.Range("D29").Value = "'2.014"
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").Value = "'1.400"
$ws.Range("E30").Value = "  +3.44%  "
$ws.Range("E31").Value = "  -1.19%  "
$ws.Range("D32").Value = "'4.250"
$ws.Range("E32").Value = "  -1.18%  "
$ws.Range("D33").Value = "'4.059"
$ws.Range("E33").Value = "  -0.39%  "
$ws.Range("D34").Value = "'0.05208"
$ws.Range("E34").Value = "  +2.08%  "
$ws.Range("D35").Value = "'1.244"
$ws.Range("E35").Value = "  -2.58%  "
$ws.Range("D36").Value = "'0.7228"
$ws.Range("E36").Value = "  -1.70%  "
$ws.Range("D37").Value = "'2.714"
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("D38").Value = "'0.01913"
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("D39").Value = "'2.769"
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("D40").Value = "'6.123"
$ws.Range("E40").Value = "  -2.67%  "
$ws.Range("D41").Value = "'0.4388"
$ws.Range("E41").Value = "  -1.56%  "
$ws.Range("D42").Value = "'71.42"
$ws.Range("E42").Value = "  -3.72%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").Value = "'0.8284"
$ws.Range("E44").Value = "  -1.07%  "
$ws.Range("E45").Value = "  -2.81%  "
$ws.Range("D46").Value = "'7.575"
$ws.Range("E46").Value = "  -0.95%  "
$ws.Range("D47").Value = "'99.71"
$ws.Range("E47").Value = "  -1.01%  "
$ws.Range("D48").Value = "'9.732"
$ws.Range("E48").Value = "  -0.46%  "
$ws.Range("D49").Value = "2.046.10"
$ws.Range("E49").Value = "  +1.06%  "
$ws.Range("E50").Value = "  -2.64%  "
$ws.Range("D51").Value = "'0.05926"
$ws.Range("E51").Value = "  -0.27%  "
